$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TFEC")

# Set the E values to 0 for rows that were previously blank
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("E13").Value = 0

# Update the view: scroll back to A1 (remove topLeftCell="C1") and change selection to E14
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E14").Select()
